$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1424363452868818
$ws.Cells.Item(2, 4).Value = 0.1322214425083246
$ws.Cells.Item(2, 5).Value = 0.1546383889056173
$ws.Cells.Item(2, 6).Value = 1.490140051058901
$ws.Cells.Item(2, 7).Value = 1.090360206271924
$ws.Cells.Item(2, 8).Value = 1.031994992964741
$ws.Cells.Item(2, 10).Value = 0.2170802489906833
$ws.Cells.Item(2, 13).Value = 5.67559464804296
$ws.Cells.Item(3, 2).Value = 0.1329334815845868
$ws.Cells.Item(3, 4).Value = 0.1275016298945104
$ws.Cells.Item(3, 5).Value = 0.1475026103910437
$ws.Cells.Item(3, 6).Value = 1.520208920582704
$ws.Cells.Item(3, 7).Value = 1.09302037860131
$ws.Cells.Item(3, 8).Value = 1.0467347102703
$ws.Cells.Item(3, 10).Value = 0.2049804097438965
$ws.Cells.Item(3, 13).Value = 5.00614054727555
$ws.Cells.Item(4, 2).Value = 0.1271690837671571
$ws.Cells.Item(4, 4).Value = 0.1246064820111883
$ws.Cells.Item(4, 5).Value = 0.1432093557517788
$ws.Cells.Item(4, 6).Value = 1.541261295821883
$ws.Cells.Item(4, 7).Value = 1.096835282291096
$ws.Cells.Item(4, 8).Value = 1.057228502865712
$ws.Cells.Item(4, 10).Value = 0.1977467703405296
$ws.Cells.Item(4, 13).Value = 4.59358577444857
$ws.Cells.Item(5, 2).Value = 0.1248378441387388
$ws.Cells.Item(5, 4).Value = 0.1234274454015747
$ws.Cells.Item(5, 5).Value = 0.1414817835958004
$ws.Cells.Item(5, 6).Value = 1.550482623319141
$ws.Cells.Item(5, 7).Value = 1.098928785314314
$ws.Cells.Item(5, 8).Value = 1.061863467395725
$ws.Cells.Item(5, 10).Value = 0.1948473396674757
$ws.Cells.Item(5, 13).Value = 4.425080118249809
$ws.Cells.Item(6, 2).Value = 0.1244518215622179
$ws.Cells.Item(6, 4).Value = 0.1232317146358213
$ws.Cells.Item(6, 5).Value = 0.1411962428126472
$ws.Cells.Item(6, 6).Value = 1.552052351670795
$ws.Cells.Item(6, 7).Value = 1.099308677203595
$ws.Cells.Item(6, 8).Value = 1.062654635977211
$ws.Cells.Item(6, 10).Value = 0.1943687863774102
$ws.Cells.Item(6, 13).Value = 4.397076370472035
$ws.Cells.Item(7, 2).Value = 0.1271375716534919
$ws.Cells.Item(7, 4).Value = 0.1245905779687035
$ws.Cells.Item(7, 5).Value = 0.1431859684221308
$ws.Cells.Item(7, 6).Value = 1.541383068918663
$ws.Cells.Item(7, 7).Value = 1.096861346716523
$ws.Cells.Item(7, 8).Value = 1.057289565105989
$ws.Cells.Item(7, 10).Value = 0.1977074730053801
$ws.Cells.Item(7, 13).Value = 4.591314820678463
$ws.Cells.Item(8, 2).Value = 0.1391451999285778
$ws.Cells.Item(8, 4).Value = 0.1305934870752168
$ws.Cells.Item(8, 5).Value = 0.1521595470580763
$ws.Cells.Item(8, 6).Value = 1.49996483536593
$ws.Cells.Item(8, 7).Value = 1.090818800388803
$ws.Cells.Item(8, 8).Value = 1.036775018318536
$ws.Cells.Item(8, 10).Value = 0.2128670531204051
$ws.Cells.Item(8, 13).Value = 5.445074737642614
$ws.Cells.Item(9, 2).Value = 0.1632480499484785
$ws.Cells.Item(9, 4).Value = 0.1423865386132945
$ws.Cells.Item(9, 5).Value = 0.170467608923218
$ws.Cells.Item(9, 6).Value = 1.439700264383077
$ws.Cells.Item(9, 7).Value = 1.09672006846867
$ws.Cells.Item(9, 8).Value = 1.008197053926551
$ws.Cells.Item(9, 10).Value = 0.2441923344177042
$ws.Cells.Item(9, 13).Value = 7.107812784689884
$ws.Cells.Item(10, 2).Value = 0.1812941243713198
$ws.Cells.Item(10, 4).Value = 0.1510637264713921
$ws.Cells.Item(10, 5).Value = 0.1843715597602653
$ws.Cells.Item(10, 6).Value = 1.408766259248338
$ws.Cells.Item(10, 7).Value = 1.112495880434437
$ws.Cells.Item(10, 8).Value = 0.9945793962921528
$ws.Cells.Item(10, 10).Value = 0.2682497608058156
$ws.Cells.Item(10, 13).Value = 8.32330563426251
$ws.Cells.Item(11, 2).Value = 0.1895770648108339
$ws.Cells.Item(11, 4).Value = 0.1550141381199808
$ws.Cells.Item(11, 5).Value = 0.1908001737484639
$ws.Cells.Item(11, 6).Value = 1.39771570467552
$ws.Cells.Item(11, 7).Value = 1.122293471826396
$ws.Cells.Item(11, 8).Value = 0.9900474581172034
$ws.Cells.Item(11, 10).Value = 0.2794370680254445
$ws.Cells.Item(11, 13).Value = 8.875163668919413
$ws.Cells.Item(12, 2).Value = 0.1927241562447648
$ws.Cells.Item(12, 4).Value = 0.1565105048719886
$ws.Cells.Item(12, 5).Value = 0.1932498457573075
$ws.Cells.Item(12, 6).Value = 1.393976157277237
$ws.Cells.Item(12, 7).Value = 1.126391899291946
$ws.Cells.Item(12, 8).Value = 0.9885755608451916
$ws.Cells.Item(12, 10).Value = 0.283709830439733
$ws.Cells.Item(12, 13).Value = 9.084001502869739
$ws.Cells.Item(13, 2).Value = 0.1920459069752241
$ws.Cells.Item(13, 4).Value = 0.1561882165137263
$ws.Cells.Item(13, 5).Value = 0.1927215781670455
$ws.Cells.Item(13, 6).Value = 1.394761570965599
$ws.Cells.Item(13, 7).Value = 1.125491776441095
$ws.Cells.Item(13, 8).Value = 0.9888816141807126
$ws.Cells.Item(13, 10).Value = 0.2827879743498727
$ws.Cells.Item(13, 13).Value = 9.039030442798605
$ws.Cells.Item(14, 2).Value = 0.1898357674569269
$ws.Cells.Item(14, 4).Value = 0.1551372365336618
$ws.Cells.Item(14, 5).Value = 0.1910014008235734
$ws.Cells.Item(14, 6).Value = 1.397399068923363
$ws.Cells.Item(14, 7).Value = 1.122622803939493
$ws.Cells.Item(14, 8).Value = 0.9899214393028046
$ws.Cells.Item(14, 10).Value = 0.2797878537076173
$ws.Cells.Item(14, 13).Value = 8.89234756865477
$ws.Cells.Item(15, 2).Value = 0.1884833612855772
$ws.Cells.Item(15, 4).Value = 0.1544935367820131
$ws.Cells.Item(15, 5).Value = 0.1899497466911555
$ws.Cells.Item(15, 6).Value = 1.399072889660133
$ws.Cells.Item(15, 7).Value = 1.120916384339097
$ws.Cells.Item(15, 8).Value = 0.9905903250570702
$ws.Cells.Item(15, 10).Value = 0.2779549717692333
$ws.Cells.Item(15, 13).Value = 8.802482492939021
$ws.Cells.Item(16, 2).Value = 0.1807542910400741
$ws.Cells.Item(16, 4).Value = 0.1508056189182696
$ws.Cells.Item(16, 5).Value = 0.1839535529196326
$ws.Cells.Item(16, 6).Value = 1.409550137819124
$ws.Cells.Item(16, 7).Value = 1.111909295614424
$ws.Cells.Item(16, 8).Value = 0.9949094659522189
$ws.Cells.Item(16, 10).Value = 0.267523657592136
$ws.Cells.Item(16, 13).Value = 8.287220085190484
$ws.Cells.Item(17, 2).Value = 0.1760315756076665
$ws.Cells.Item(17, 4).Value = 0.1485439849378309
$ws.Cells.Item(17, 5).Value = 0.1803018944789301
$ws.Cells.Item(17, 6).Value = 1.416758919699731
$ws.Cells.Item(17, 7).Value = 1.107063089581203
$ws.Cells.Item(17, 8).Value = 0.9979886269167935
$ws.Cells.Item(17, 10).Value = 0.2611876835722171
$ws.Cells.Item(17, 13).Value = 7.970858523457878
$ws.Cells.Item(18, 2).Value = 0.1733221348472114
$ws.Cells.Item(18, 4).Value = 0.1472434466896999
$ws.Cells.Item(18, 5).Value = 0.1782112863064782
$ws.Cells.Item(18, 6).Value = 1.421188971678944
$ws.Cells.Item(18, 7).Value = 1.104521695698651
$ws.Cells.Item(18, 8).Value = 0.9999158761484637
$ws.Cells.Item(18, 10).Value = 0.2575662281486757
$ws.Cells.Item(18, 13).Value = 7.788793412625182
$ws.Cells.Item(19, 2).Value = 0.1724059591305007
$ws.Cells.Item(19, 4).Value = 0.1468031578528297
$ws.Cells.Item(19, 5).Value = 0.1775051019011826
$ws.Cells.Item(19, 6).Value = 1.422737326935845
$ws.Cells.Item(19, 7).Value = 1.103703133731983
$ws.Cells.Item(19, 8).Value = 1.000595081696048
$ws.Cells.Item(19, 10).Value = 0.2563439463796726
$ws.Cells.Item(19, 13).Value = 7.727131266826234
$ws.Cells.Item(20, 2).Value = 0.1765335990014307
$ws.Cells.Item(20, 4).Value = 0.1487847094684724
$ws.Cells.Item(20, 5).Value = 0.1806896100882582
$ws.Cells.Item(20, 6).Value = 1.415962089481596
$ws.Cells.Item(20, 7).Value = 1.107553433099866
$ws.Cells.Item(20, 8).Value = 0.9976446422578817
$ws.Cells.Item(20, 10).Value = 0.2618597864840524
$ws.Cells.Item(20, 13).Value = 8.004546230589312
$ws.Cells.Item(21, 2).Value = 0.1904846540224838
$ws.Cells.Item(21, 4).Value = 0.1554459232153818
$ws.Cells.Item(21, 5).Value = 0.1915062401917709
$ws.Cells.Item(21, 6).Value = 1.396612207764392
$ws.Cells.Item(21, 7).Value = 1.123454860691055
$ws.Cells.Item(21, 8).Value = 0.9896093470985079
$ws.Cells.Item(21, 10).Value = 0.2806680634606522
$ws.Cells.Item(21, 13).Value = 8.935435517460235
$ws.Cells.Item(22, 2).Value = 0.1996637724413262
$ws.Cells.Item(22, 4).Value = 0.1598019242019291
$ws.Cells.Item(22, 5).Value = 0.1986648937523583
$ws.Cells.Item(22, 6).Value = 1.386565141142654
$ws.Cells.Item(22, 7).Value = 1.136116079260006
$ws.Cells.Item(22, 8).Value = 0.985783952423759
$ws.Cells.Item(22, 10).Value = 0.2931730564125701
$ws.Cells.Item(22, 13).Value = 9.54302484357828
$ws.Cells.Item(23, 2).Value = 0.194759118769241
$ws.Cells.Item(23, 4).Value = 0.1574768182892115
$ws.Cells.Item(23, 5).Value = 0.1948358703229189
$ws.Cells.Item(23, 6).Value = 1.391686072842219
$ws.Cells.Item(23, 7).Value = 1.129147138521091
$ws.Cells.Item(23, 8).Value = 0.9876934194796831
$ws.Cells.Item(23, 10).Value = 0.2864789723007135
$ws.Cells.Item(23, 13).Value = 9.218810567473099
$ws.Cells.Item(24, 2).Value = 0.1763066164297129
$ws.Cells.Item(24, 4).Value = 0.1486758788245055
$ws.Cells.Item(24, 5).Value = 0.1805142965461926
$ws.Cells.Item(24, 6).Value = 1.416321447555944
$ws.Cells.Item(24, 7).Value = 1.107330987298013
$ws.Cells.Item(24, 8).Value = 0.9977996690942632
$ws.Cells.Item(24, 10).Value = 0.2615558628870787
$ws.Cells.Item(24, 13).Value = 7.989316594468505
$ws.Cells.Item(25, 2).Value = 0.156668218767976
$ws.Cells.Item(25, 4).Value = 0.139193962424045
$ws.Cells.Item(25, 5).Value = 0.1654367777617338
$ws.Cells.Item(25, 6).Value = 1.453701928048076
$ws.Cells.Item(25, 7).Value = 1.093159998725127
$ws.Cells.Item(25, 8).Value = 1.014652629159997
$ws.Cells.Item(25, 10).Value = 0.2355401676461213
$ws.Cells.Item(25, 13).Value = 6.659142366164247
